# feat: add 2022-Q3 data
#
# Inserts a new "2022-Q3" fund-holdings worksheet between the "总计"
# (totals) summary sheet and the existing "2021-Q1" worksheet, and updates
# the totals sheet with a new summary row for the added quarter.

$wb = $excel.ActiveWorkbook

# Writes $text into $targetAddr as a genuine text value (no numeric/leading
# -zero coercion) without leaving a lasting style change on the target
# cell: stage the text (via a leading apostrophe) in an unused scratch
# cell, copy only its VALUE (xlPasteValues = -4163) onto the target, then
# wipe the scratch cell completely (Clear = contents + formats) so nothing
# about it lingers in the saved workbook.
function Set-TextValue($ws, $targetAddr, $text) {
    $scratch = $ws.Range("ZZ1")
    $scratch.Value = "'" + $text
    $scratch.Copy()
    $ws.Range($targetAddr).PasteSpecial(-4163)
    $scratch.Clear()
}

$totals = $wb.Worksheets.Item("总计")
$quarter = $wb.Worksheets.Item("2021-Q1")

# --- 1. Preserve the existing "2021-Q1" sheet ---------------------------
# Duplicate it (with all of its data/formatting intact) and place the copy
# right after the original; the copy keeps the "2021-Q1" name while the
# original slot turns into the new "2022-Q3" sheet below.
$quarter.Copy($null, $quarter)
$quarterCopy = $quarter.Next
$quarterCopy.Name = "2021-Q1-restored"

# --- 2. Turn the original sheet (now at position 2) into "2022-Q3" ------
$quarter.Name = "2022-Q3"

# Match the header/first-column format to the one used by the "总计"
# sheet (the style newly-appended quarter sheets use).
$totals.Range("B1").Copy()
$quarter.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats
$totals.Range("A2").Copy()
$quarter.Range("A2").PasteSpecial(-4122)      # xlPasteFormats

$quarter.Range("B1").Value = "基金代码"
$quarter.Range("C1").Value = "基金名称"
$quarter.Range("D1").Value = "基金规模"
$quarter.Range("E1").Value = "股票总仓位"
$quarter.Range("F1").Value = "仓位占比"
$quarter.Range("G1").Value = "持有市值(亿元)"
$quarter.Range("H1").Value = "仓位排名"

$quarter.Range("A2").Value = 0
Set-TextValue $quarter "B2" "159628"
$quarter.Range("C2").Value = "万家国证2000ETF"
Set-TextValue $quarter "D2" "2.90"
Set-TextValue $quarter "E2" "97.72"
Set-TextValue $quarter "F2" "0.45"
Set-TextValue $quarter "G2" "0.0130"
$quarter.Range("H2").Value = 8

# --- 3. Rename the restored copy back to "2021-Q1" ----------------------
$quarterCopy.Name = "2021-Q1"

# --- 4. Update the "总计" (totals) sheet ---------------------------------
# Row 2 used to describe "2021-Q1"; it now describes the new "2022-Q3"
# quarter, and a fresh row 3 is appended with the original "2021-Q1"
# totals (same formatting as row 2's first column).
$totals.Range("A2").Copy()
$totals.Range("A3").PasteSpecial(-4122)       # xlPasteFormats
$totals.Range("A3").Value = 1
$totals.Range("B3").Value = "2021-Q1"
$totals.Range("C3").Value = 1
$totals.Range("D3").Value = 0

$totals.Range("B2").Value = "2022-Q3"
$totals.Range("D2").Value = 0.01
